$wb = $excel.ActiveWorkbook

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2333.6458
$ws.Range("I58").Value = 332.5
$ws.Range("J58").Value = 2619.524
$ws.Range("K58").Value = 997.5
$ws.Range("L58").Value = 7858.572
$ws.Range("M58").Value = -847.5
$ws.Range("N58").Value = -8158.572

# ALC row 82
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 7625
$ws.Range("J82").Value = 8000
$ws.Range("L82").Value = 24000
$ws.Range("N82").Value = -24812

# ALC row 85
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 7625
$ws.Range("J85").Value = 8000
$ws.Range("L85").Value = 24000
$ws.Range("M85").Value = -20346
$ws.Range("N85").Value = -26808

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1009.6667
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1014.5
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 3043.5
$ws.Range("M96").Value = -1627
$ws.Range("N96").Value = -5789.5

# ALC row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1855.4166
$ws.Range("I115").Value = 421.66666
$ws.Range("K115").Value = 1264.99998
$ws.Range("M115").Value = 302.0000199999999

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1482.8667
$ws.Range("I135").Value = 1445.9286
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 13013.3574
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -10478.3574
$ws.Range("N135").Value = -23070

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3424.1072
$ws.Range("I137").Value = 3180.5557
$ws.Range("J137").Value = 10000
$ws.Range("K137").Value = 9541.667099999999
$ws.Range("L137").Value = 30000
$ws.Range("M137").Value = -6991.667099999999
$ws.Range("N137").Value = -35100

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 924.1739
$ws.Range("I74").Value = 873.73334
$ws.Range("K74").Value = 873.73334
$ws.Range("M74").Value = 0.2666600000000017

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 924.1739
$ws.Range("I77").Value = 873.73334
$ws.Range("K77").Value = 4368.6667
$ws.Range("M77").Value = -0.6666999999997643

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2627
$ws.Range("I122").Value = 1826.5454
$ws.Range("K122").Value = 5479.6362
$ws.Range("M122").Value = -3029.6362

# ARM row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 29623.75
$ws.Range("J137").Value = 29623.75
$ws.Range("L137").Value = 29623.75
$ws.Range("N137").Value = -39823.75

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2864.0557
$ws.Range("I134").Value = 1809.6154
$ws.Range("K134").Value = 5428.8462
$ws.Range("M134").Value = -2893.8462

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4067.0408
$ws.Range("I31").Value = 3000.484
$ws.Range("J31").Value = 5903.8887
$ws.Range("K31").Value = 3000.484
$ws.Range("L31").Value = 5903.8887
$ws.Range("M31").Value = -2705.484
$ws.Range("N31").Value = -6493.8887

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4067.0408
$ws.Range("I34").Value = 3000.484
$ws.Range("J34").Value = 5903.8887
$ws.Range("K34").Value = 3000.484
$ws.Range("L34").Value = 5903.8887
$ws.Range("M34").Value = -2798.484
$ws.Range("N34").Value = -6307.8887

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3174.1
$ws.Range("I122").Value = 2731.923
$ws.Range("J122").Value = 3995.2856
$ws.Range("K122").Value = 8195.769
$ws.Range("L122").Value = 11985.8568
$ws.Range("M122").Value = -5745.769
$ws.Range("N122").Value = -16885.8568

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2123.2
$ws.Range("J75").Value = 2599
$ws.Range("L75").Value = 7797
$ws.Range("N75").Value = -9793

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2123.2
$ws.Range("J78").Value = 2599
$ws.Range("L78").Value = 23391
$ws.Range("N78").Value = -33375

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4480.6
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4480.6
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 13441.8
$ws.Range("N80").Value = -15313.8
$ws.Range("M80").ClearContents()

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4480.6
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4480.6
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 40325.4
$ws.Range("N83").Value = -49685.4
$ws.Range("M83").ClearContents()

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11389.866
$ws.Range("J87").Value = 15644.444
$ws.Range("L87").Value = 46933.33199999999
$ws.Range("N87").Value = -49429.33199999999

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 11389.866
$ws.Range("J90").Value = 15644.444
$ws.Range("L90").Value = 140799.996
$ws.Range("N90").Value = -153279.996

# CUL row 120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 18970.334

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 58503.723
$ws.Range("I102").Value = 2582.2307
$ws.Range("J102").Value = 203899.6
$ws.Range("K102").Value = 2582.2307
$ws.Range("L102").Value = 203899.6
$ws.Range("M102").Value = -960.2307000000001
$ws.Range("N102").Value = -207143.6

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2681.818
$ws.Range("I7").Value = 1937.5
$ws.Range("J7").Value = 4666.6665
$ws.Range("K7").Value = 1937.5
$ws.Range("L7").Value = 4666.6665
$ws.Range("M7").Value = -1825.5
$ws.Range("N7").Value = -4890.6665

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 111114320
$ws.Range("I61").Value = 142859570
$ws.Range("J61").Value = 5950.5
$ws.Range("K61").Value = 142859570
$ws.Range("L61").Value = 5950.5
$ws.Range("M61").Value = -142859368
$ws.Range("N61").Value = -6354.5

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3340.7058
$ws.Range("I93").Value = 2485.1428
$ws.Range("J93").Value = 7333.3335
$ws.Range("K93").Value = 2485.1428
$ws.Range("L93").Value = 7333.3335
$ws.Range("M93").Value = -1237.1428
$ws.Range("N93").Value = -9829.333500000001

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 111114320
$ws.Range("I113").Value = 142859570
$ws.Range("J113").Value = 5950.5
$ws.Range("K113").Value = 142859570
$ws.Range("L113").Value = 5950.5
$ws.Range("M113").Value = -142857400
$ws.Range("N113").Value = -10290.5

# LTW row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2681.818
$ws.Range("I126").Value = 1937.5
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 5812.5
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -3342.5
$ws.Range("N126").Value = -18939.9995

# WVR row 68
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 39600
$ws.Range("J68").Value = 39600
$ws.Range("L68").Value = 39600
$ws.Range("N68").Value = -41222

# WVR row 71
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 39600
$ws.Range("J71").Value = 39600
$ws.Range("L71").Value = 118800
$ws.Range("N71").Value = -126912

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13928.303
$ws.Range("I132").Value = 2197
$ws.Range("K132").Value = 6591
$ws.Range("M132").Value = -4061
